$wb = $excel.ActiveWorkbook

# "_set_PRODUCT_DATA" sheet: remove the second data row
# (pd_Names="unit energy use", pd_Category="energy_use") which is no
# longer needed now that "unit energy use, initial" / "energy_use_0"
# covers it. Remaining rows shift up automatically.
$wsProductData = $wb.Worksheets.Item("_set_PRODUCT_DATA")
$wsProductData.Rows.Item(2).Delete()

# Make "_set_PRODUCT_DATA" the active sheet/tab and move the selection
# further down the (now shorter) sheet.
$wsProductData.Activate()
$wsProductData.Range("A20").Select()
